$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old row 2 data (A2:F2)
$ws.Range("A2:F2").ClearContents()

# Write the new row 11 data (A11:D11)
$ws.Range("A11").Value = 10.0
$ws.Range("B11").Value = 23.859
$ws.Range("C11").Value = 46.0
$ws.Range("D11").Value = 1.6742454441366514
